$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Fix 1: correct the typo "#AdmissionReqiurements" -> "#AdmissionRequirements"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("#AdmissionReqiurements", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "#AdmissionRequirements", 2)

# ---------------------------------------------------------------------------
# Fix 2: the merge field "#DegreeRequiredCredits" used to live in a single
# run; the placeholder hash and the field name now need to be split into two
# separate runs (the "#" keeps the field's original run, the field name
# becomes its own run) so that requiredCredits can be re-pointed at the new
# Specialization-level field while leaving the leading "#" alone.
# Only the FIRST occurrence (the one carrying <w:lang w:val="en-US"/>) is
# affected - a second, unrelated "#DegreeRequiredCredits" later in the
# document must be left untouched.
# ---------------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("#DegreeRequiredCredits")

if ($found) {
    $matchStart = $target.Start
    $matchEnd = $target.End

    # Split the run after the leading "#" by toggling a character property on
    # the single "#" character - this forces the run to break in two while
    # each half keeps the original run formatting. Bold is left ON for now
    # so the two halves stay distinct (same-formatting adjacent runs can get
    # silently re-merged by the next step otherwise).
    $hashChar = $d.Range($matchStart, $matchStart + 1)
    $hashChar.Bold = 1

    # Re-stamp the trailing "DegreeRequiredCredits" run (a plain Find/Replace
    # of the exact text back onto itself re-mints the run as a fresh one, the
    # way it appears after being retargeted to the new field).
    $fieldNameRange = $d.Range($matchStart + 1, $matchEnd)
    $fieldNameText = $fieldNameRange.Text
    $fieldNameRange.Find.Execute($fieldNameText, $false, $false, $false, $false, $false, `
                                  $true, 0, $false, $fieldNameText, 1)

    # Now restore the leading "#" back to non-bold, leaving it as its own run.
    $hashChar2 = $d.Range($matchStart, $matchStart + 1)
    $hashChar2.Bold = 0
}

Write-Host "Done"
